# Updates cryptos list cell values (Price + Volume(1h) columns) to match
# the latest scrape. Values are stored as plain text in the sheet (same
# as the original inline-string cells), so numeric-looking Price values
# get the cell format set to Text first -- otherwise Excel's COM layer
# auto-converts a plain numeric string into a real number on assignment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.700.79'
$ws.Range('E2').Value = '  +0.77%  '
$ws.Range('D3').Value = '2.444.36'
$ws.Range('E3').Value = '  +1.05%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.92'
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.26'
$ws.Range('E6').Value = '  +2.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  +0.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.113'
$ws.Range('E9').Value = '  +3.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.28'
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.354'
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000187'
$ws.Range('E13').Value = '  +6.95%  '
$ws.Range('E14').Value = '  +4.40%  '
$ws.Range('D16').Value = '62.478.96'
$ws.Range('E16').Value = '  +0.48%  '
$ws.Range('D17').Value = '2.444.02'
$ws.Range('E17').Value = '  +1.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.33'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('E19').Value = '  +1.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '325.52'
$ws.Range('E20').Value = '  +0.42%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.46'
$ws.Range('E23').Value = '  +2.67%  '
$ws.Range('E24').Value = '  +2.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.86'
$ws.Range('E25').Value = '  -2.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '570.65'
$ws.Range('E26').Value = '  -1.41%  '
$ws.Range('D27').Value = '0.0₃0987'
$ws.Range('E27').Value = '  +3.60%  '
$ws.Range('D28').Value = '2.563.51'
$ws.Range('E28').Value = '  +1.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('E30').Value = '  +2.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.47'
$ws.Range('E31').Value = '  +1.68%  '
$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('E33').Value = '  +0.32%  '
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.90'
$ws.Range('E35').Value = '  +3.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.997'
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('E37').Value = '  +0.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.58'
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.82'
$ws.Range('E39').Value = '  +0.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '150.60'
$ws.Range('E40').Value = '  -1.43%  '
$ws.Range('E41').Value = '  +0.61%  '
$ws.Range('E42').Value = '  +0.63%  '
$ws.Range('E43').Value = '  +5.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '149.65'
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.71'
$ws.Range('E45').Value = '  +1.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0539'
$ws.Range('E46').Value = '  +0.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '20.49'
$ws.Range('E47').Value = '  +1.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.602'
$ws.Range('E48').Value = '  +1.13%  '
$ws.Range('E49').Value = '  +1.48%  '
$ws.Range('E50').Value = '  +1.49%  '
$ws.Range('E51').Value = '  +0.69%  '
